$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing row (613) down to the new rows (614-627)
$ws.Range("A613:V613").Copy() | Out-Null
$ws.Range("A614:V627").PasteSpecial(-4122) | Out-Null

# Populate the new rows with data
# Row 614
$ws.Cells.Item(614, 1).Value = "Entrainement"
$ws.Cells.Item(614, 2).Value = 45934
$ws.Cells.Item(614, 3).Value = "Global"
$ws.Cells.Item(614, 4).Value = "M"
$ws.Cells.Item(614, 5).Value = "Sofiane Belle"
$ws.Cells.Item(614, 6).Value = "left forward"
$ws.Cells.Item(614, 7).Value = "01:25:25"
$ws.Cells.Item(614, 8).Value = 9.08
$ws.Cells.Item(614, 9).Value = 2.17
$ws.Cells.Item(614, 10).Value = 6.89
$ws.Cells.Item(614, 11).Value = 1.28
$ws.Cells.Item(614, 12).Value = 0.75
$ws.Cells.Item(614, 13).Value = 0.15
$ws.Cells.Item(614, 14).Value = 0.02
$ws.Cells.Item(614, 15).Value = 10
$ws.Cells.Item(614, 16).Value = 6.36
$ws.Cells.Item(614, 17).Value = 31.13
$ws.Cells.Item(614, 18).Value = 4.28
$ws.Cells.Item(614, 19).Value = 30
$ws.Cells.Item(614, 20).Value = 6
$ws.Cells.Item(614, 21).Value = 36
$ws.Cells.Item(614, 22).Value = 15

# Row 615
$ws.Cells.Item(615, 1).Value = "Entrainement"
$ws.Cells.Item(615, 2).Value = 45934
$ws.Cells.Item(615, 3).Value = "Global"
$ws.Cells.Item(615, 4).Value = "M"
$ws.Cells.Item(615, 5).Value = "Karahali Souaré"
$ws.Cells.Item(615, 6).Value = "right forward"
$ws.Cells.Item(615, 7).Value = "00:56:04"
$ws.Cells.Item(615, 8).Value = 5.68
$ws.Cells.Item(615, 9).Value = 1.25
$ws.Cells.Item(615, 10).Value = 4.42
$ws.Cells.Item(615, 11).Value = 0.74
$ws.Cells.Item(615, 12).Value = 0.32
$ws.Cells.Item(615, 13).Value = 0.16
$ws.Cells.Item(615, 14).Value = 0.05
$ws.Cells.Item(615, 15).Value = 12
$ws.Cells.Item(615, 16).Value = 6.04
$ws.Cells.Item(615, 17).Value = 32.33
$ws.Cells.Item(615, 18).Value = 4.73
$ws.Cells.Item(615, 19).Value = 37
$ws.Cells.Item(615, 20).Value = 10
$ws.Cells.Item(615, 21).Value = 26
$ws.Cells.Item(615, 22).Value = 19

# Row 616
$ws.Cells.Item(616, 1).Value = "Entrainement"
$ws.Cells.Item(616, 2).Value = 45934
$ws.Cells.Item(616, 3).Value = "Global"
$ws.Cells.Item(616, 4).Value = "M"
$ws.Cells.Item(616, 5).Value = "Amir Etien"
$ws.Cells.Item(616, 6).Value = "right forward"
$ws.Cells.Item(616, 7).Value = "01:25:33"
$ws.Cells.Item(616, 8).Value = 8.94
$ws.Cells.Item(616, 9).Value = 1.93
$ws.Cells.Item(616, 10).Value = 6.98
$ws.Cells.Item(616, 11).Value = 1.01
$ws.Cells.Item(616, 12).Value = 0.61
$ws.Cells.Item(616, 13).Value = 0.28
$ws.Cells.Item(616, 14).Value = 0.06
$ws.Cells.Item(616, 15).Value = 16
$ws.Cells.Item(616, 16).Value = 6.18
$ws.Cells.Item(616, 17).Value = 33.35
$ws.Cells.Item(616, 18).Value = 4.58
$ws.Cells.Item(616, 19).Value = 42
$ws.Cells.Item(616, 20).Value = 11
$ws.Cells.Item(616, 21).Value = 25
$ws.Cells.Item(616, 22).Value = 12

# Row 617
$ws.Cells.Item(617, 1).Value = "Entrainement"
$ws.Cells.Item(617, 2).Value = 45934
$ws.Cells.Item(617, 3).Value = "Global"
$ws.Cells.Item(617, 4).Value = "M"
$ws.Cells.Item(617, 5).Value = "Ilyes Boughanmi"
$ws.Cells.Item(617, 6).Value = "center forward"
$ws.Cells.Item(617, 7).Value = "00:13:45"
$ws.Cells.Item(617, 8).Value = 1.39
$ws.Cells.Item(617, 9).Value = 0.26
$ws.Cells.Item(617, 10).Value = 1.13
$ws.Cells.Item(617, 11).Value = 0.15
$ws.Cells.Item(617, 12).Value = 0.09
$ws.Cells.Item(617, 13).Value = 0.02
$ws.Cells.Item(617, 14).Value = 0
$ws.Cells.Item(617, 15).Value = 4
$ws.Cells.Item(617, 16).Value = 6.07
$ws.Cells.Item(617, 17).Value = 25.61
$ws.Cells.Item(617, 18).Value = 3.38
$ws.Cells.Item(617, 19).Value = 5
$ws.Cells.Item(617, 20).Value = 0
$ws.Cells.Item(617, 21).Value = 8
$ws.Cells.Item(617, 22).Value = 1

# Row 618
$ws.Cells.Item(618, 1).Value = "Entrainement"
$ws.Cells.Item(618, 2).Value = 45934
$ws.Cells.Item(618, 3).Value = "Global"
$ws.Cells.Item(618, 4).Value = "M"
$ws.Cells.Item(618, 5).Value = "Mattheo Haon"
$ws.Cells.Item(618, 6).Value = "right back"
$ws.Cells.Item(618, 7).Value = "01:45:35"
$ws.Cells.Item(618, 8).Value = 11.67
$ws.Cells.Item(618, 9).Value = 2.3
$ws.Cells.Item(618, 10).Value = 9.34
$ws.Cells.Item(618, 11).Value = 1.56
$ws.Cells.Item(618, 12).Value = 0.61
$ws.Cells.Item(618, 13).Value = 0.14
$ws.Cells.Item(618, 14).Value = 0.01
$ws.Cells.Item(618, 15).Value = 12
$ws.Cells.Item(618, 16).Value = 6.55
$ws.Cells.Item(618, 17).Value = 31.06
$ws.Cells.Item(618, 18).Value = 4.81
$ws.Cells.Item(618, 19).Value = 57
$ws.Cells.Item(618, 20).Value = 9
$ws.Cells.Item(618, 21).Value = 38
$ws.Cells.Item(618, 22).Value = 13

# Row 619
$ws.Cells.Item(619, 1).Value = "Entrainement"
$ws.Cells.Item(619, 2).Value = 45934
$ws.Cells.Item(619, 3).Value = "Global"
$ws.Cells.Item(619, 4).Value = "M"
$ws.Cells.Item(619, 5).Value = "Karim Belmahi"
$ws.Cells.Item(619, 6).Value = "left forward"
$ws.Cells.Item(619, 7).Value = "00:14:25"
$ws.Cells.Item(619, 8).Value = 1.62
$ws.Cells.Item(619, 9).Value = 0.3
$ws.Cells.Item(619, 10).Value = 1.31
$ws.Cells.Item(619, 11).Value = 0.17
$ws.Cells.Item(619, 12).Value = 0.1
$ws.Cells.Item(619, 13).Value = 0.03
$ws.Cells.Item(619, 14).Value = 0
$ws.Cells.Item(619, 15).Value = 2
$ws.Cells.Item(619, 16).Value = 6.69
$ws.Cells.Item(619, 17).Value = 26.38
$ws.Cells.Item(619, 18).Value = 4.34
$ws.Cells.Item(619, 19).Value = 17
$ws.Cells.Item(619, 20).Value = 2
$ws.Cells.Item(619, 21).Value = 9
$ws.Cells.Item(619, 22).Value = 2

# Row 620
$ws.Cells.Item(620, 1).Value = "Entrainement"
$ws.Cells.Item(620, 2).Value = 45934
$ws.Cells.Item(620, 3).Value = "Global"
$ws.Cells.Item(620, 4).Value = "M"
$ws.Cells.Item(620, 5).Value = "Naim Ighbane"
$ws.Cells.Item(620, 6).Value = "center back"
$ws.Cells.Item(620, 7).Value = "01:45:35"
$ws.Cells.Item(620, 8).Value = 10.58
$ws.Cells.Item(620, 9).Value = 1.66
$ws.Cells.Item(620, 10).Value = 8.9
$ws.Cells.Item(620, 11).Value = 1.12
$ws.Cells.Item(620, 12).Value = 0.47
$ws.Cells.Item(620, 13).Value = 0.08
$ws.Cells.Item(620, 14).Value = 0
$ws.Cells.Item(620, 15).Value = 7
$ws.Cells.Item(620, 16).Value = 5.89
$ws.Cells.Item(620, 17).Value = 29.6
$ws.Cells.Item(620, 18).Value = 4.63
$ws.Cells.Item(620, 19).Value = 47
$ws.Cells.Item(620, 20).Value = 4
$ws.Cells.Item(620, 21).Value = 27
$ws.Cells.Item(620, 22).Value = 12

# Row 621
$ws.Cells.Item(621, 1).Value = "Entrainement"
$ws.Cells.Item(621, 2).Value = 45934
$ws.Cells.Item(621, 3).Value = "Global"
$ws.Cells.Item(621, 4).Value = "M"
$ws.Cells.Item(621, 5).Value = "Kamal Bafounta"
$ws.Cells.Item(621, 6).Value = "center midfield"
$ws.Cells.Item(621, 7).Value = "00:49:15"
$ws.Cells.Item(621, 8).Value = 6.2
$ws.Cells.Item(621, 9).Value = 1.27
$ws.Cells.Item(621, 10).Value = 4.92
$ws.Cells.Item(621, 11).Value = 0.95
$ws.Cells.Item(621, 12).Value = 0.28
$ws.Cells.Item(621, 13).Value = 0.05
$ws.Cells.Item(621, 14).Value = 0.01
$ws.Cells.Item(621, 15).Value = 2
$ws.Cells.Item(621, 16).Value = 7.61
$ws.Cells.Item(621, 17).Value = 30.4
$ws.Cells.Item(621, 18).Value = 4.44
$ws.Cells.Item(621, 19).Value = 20
$ws.Cells.Item(621, 20).Value = 3
$ws.Cells.Item(621, 21).Value = 25
$ws.Cells.Item(621, 22).Value = 8

# Row 622
$ws.Cells.Item(622, 1).Value = "Entrainement"
$ws.Cells.Item(622, 2).Value = 45934
$ws.Cells.Item(622, 3).Value = "Global"
$ws.Cells.Item(622, 4).Value = "M"
$ws.Cells.Item(622, 5).Value = "Malik Boussaid"
$ws.Cells.Item(622, 6).Value = "right back"
$ws.Cells.Item(622, 7).Value = "01:05:01"
$ws.Cells.Item(622, 8).Value = 6.73
$ws.Cells.Item(622, 9).Value = 1.42
$ws.Cells.Item(622, 10).Value = 5.29
$ws.Cells.Item(622, 11).Value = 0.96
$ws.Cells.Item(622, 12).Value = 0.38
$ws.Cells.Item(622, 13).Value = 0.1
$ws.Cells.Item(622, 14).Value = 0
$ws.Cells.Item(622, 15).Value = 7
$ws.Cells.Item(622, 16).Value = 6.1
$ws.Cells.Item(622, 17).Value = 28.08
$ws.Cells.Item(622, 18).Value = 4.76
$ws.Cells.Item(622, 19).Value = 24
$ws.Cells.Item(622, 20).Value = 3
$ws.Cells.Item(622, 21).Value = 27
$ws.Cells.Item(622, 22).Value = 7

# Row 623
$ws.Cells.Item(623, 1).Value = "Entrainement"
$ws.Cells.Item(623, 2).Value = 45934
$ws.Cells.Item(623, 3).Value = "Global"
$ws.Cells.Item(623, 4).Value = "M"
$ws.Cells.Item(623, 5).Value = "Ilan Ihaddadene"
$ws.Cells.Item(623, 6).Value = "center midfield"
$ws.Cells.Item(623, 7).Value = "01:45:35"
$ws.Cells.Item(623, 8).Value = 12.94
$ws.Cells.Item(623, 9).Value = 2.76
$ws.Cells.Item(623, 10).Value = 10.15
$ws.Cells.Item(623, 11).Value = 2.13
$ws.Cells.Item(623, 12).Value = 0.5
$ws.Cells.Item(623, 13).Value = 0.15
$ws.Cells.Item(623, 14).Value = 0.02
$ws.Cells.Item(623, 15).Value = 9
$ws.Cells.Item(623, 16).Value = 7.29
$ws.Cells.Item(623, 17).Value = 31.28
$ws.Cells.Item(623, 18).Value = 4.61
$ws.Cells.Item(623, 19).Value = 54
$ws.Cells.Item(623, 20).Value = 15
$ws.Cells.Item(623, 21).Value = 41
$ws.Cells.Item(623, 22).Value = 11

# Row 624
$ws.Cells.Item(624, 1).Value = "Entrainement"
$ws.Cells.Item(624, 2).Value = 45934
$ws.Cells.Item(624, 3).Value = "Global"
$ws.Cells.Item(624, 4).Value = "M"
$ws.Cells.Item(624, 5).Value = "Emmanuel Valey"
$ws.Cells.Item(624, 6).Value = "left forward"
$ws.Cells.Item(624, 7).Value = "01:12:53"
$ws.Cells.Item(624, 8).Value = 9.24
$ws.Cells.Item(624, 9).Value = 2.15
$ws.Cells.Item(624, 10).Value = 7.06
$ws.Cells.Item(624, 11).Value = 1.21
$ws.Cells.Item(624, 12).Value = 0.68
$ws.Cells.Item(624, 13).Value = 0.23
$ws.Cells.Item(624, 14).Value = 0.06
$ws.Cells.Item(624, 15).Value = 19
$ws.Cells.Item(624, 16).Value = 7.58
$ws.Cells.Item(624, 17).Value = 32.84
$ws.Cells.Item(624, 18).Value = 5.36
$ws.Cells.Item(624, 19).Value = 46
$ws.Cells.Item(624, 20).Value = 7
$ws.Cells.Item(624, 21).Value = 42
$ws.Cells.Item(624, 22).Value = 14

# Row 625
$ws.Cells.Item(625, 1).Value = "Entrainement"
$ws.Cells.Item(625, 2).Value = 45934
$ws.Cells.Item(625, 3).Value = "Global"
$ws.Cells.Item(625, 4).Value = "M"
$ws.Cells.Item(625, 5).Value = "Yoan Zouma"
$ws.Cells.Item(625, 6).Value = "center back"
$ws.Cells.Item(625, 7).Value = "01:45:19"
$ws.Cells.Item(625, 8).Value = 9.65
$ws.Cells.Item(625, 9).Value = 1.2
$ws.Cells.Item(625, 10).Value = 8.44
$ws.Cells.Item(625, 11).Value = 0.75
$ws.Cells.Item(625, 12).Value = 0.29
$ws.Cells.Item(625, 13).Value = 0.15
$ws.Cells.Item(625, 14).Value = 0.02
$ws.Cells.Item(625, 15).Value = 12
$ws.Cells.Item(625, 16).Value = 5.42
$ws.Cells.Item(625, 17).Value = 31.16
$ws.Cells.Item(625, 18).Value = 5.1
$ws.Cells.Item(625, 19).Value = 21
$ws.Cells.Item(625, 20).Value = 4
$ws.Cells.Item(625, 21).Value = 25
$ws.Cells.Item(625, 22).Value = 2

# Row 626
$ws.Cells.Item(626, 1).Value = "Entrainement"
$ws.Cells.Item(626, 2).Value = 45934
$ws.Cells.Item(626, 3).Value = "Global"
$ws.Cells.Item(626, 4).Value = "M"
$ws.Cells.Item(626, 5).Value = "Naim Dhib"
$ws.Cells.Item(626, 6).Value = "center midfield"
$ws.Cells.Item(626, 7).Value = "01:45:35"
$ws.Cells.Item(626, 8).Value = 10.97
$ws.Cells.Item(626, 9).Value = 2
$ws.Cells.Item(626, 10).Value = 8.94
$ws.Cells.Item(626, 11).Value = 1.55
$ws.Cells.Item(626, 12).Value = 0.43
$ws.Cells.Item(626, 13).Value = 0.04
$ws.Cells.Item(626, 14).Value = 0
$ws.Cells.Item(626, 15).Value = 4
$ws.Cells.Item(626, 16).Value = 6.15
$ws.Cells.Item(626, 17).Value = 28.51
$ws.Cells.Item(626, 18).Value = 5.49
$ws.Cells.Item(626, 19).Value = 45
$ws.Cells.Item(626, 20).Value = 9
$ws.Cells.Item(626, 21).Value = 45
$ws.Cells.Item(626, 22).Value = 17

# Row 627
$ws.Cells.Item(627, 1).Value = "Entrainement"
$ws.Cells.Item(627, 2).Value = 45934
$ws.Cells.Item(627, 3).Value = "Global"
$ws.Cells.Item(627, 4).Value = "M"
$ws.Cells.Item(627, 5).Value = "Levy Ndoutoume"
$ws.Cells.Item(627, 6).Value = "left back"
$ws.Cells.Item(627, 7).Value = "00:40:50"
$ws.Cells.Item(627, 8).Value = 4.26
$ws.Cells.Item(627, 9).Value = 0.85
$ws.Cells.Item(627, 10).Value = 3.4
$ws.Cells.Item(627, 11).Value = 0.38
$ws.Cells.Item(627, 12).Value = 0.31
$ws.Cells.Item(627, 13).Value = 0.17
$ws.Cells.Item(627, 14).Value = 0
$ws.Cells.Item(627, 15).Value = 11
$ws.Cells.Item(627, 16).Value = 6.23
$ws.Cells.Item(627, 17).Value = 29.74
$ws.Cells.Item(627, 18).Value = 4.55
$ws.Cells.Item(627, 19).Value = 21
$ws.Cells.Item(627, 20).Value = 5
$ws.Cells.Item(627, 21).Value = 19
$ws.Cells.Item(627, 22).Value = 9

# Update selection / view to match the final state
$ws.Range("C631").Select()
